$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeRef, $val) {
    $rng = $ws.Range($rangeRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '59.186.98'
Set-TextValue 'E2' '  -0.52%  '
Set-TextValue 'D3' '2.637.94'
Set-TextValue 'E3' '  -0.34%  '
Set-TextValue 'E4' '  +0.18%  '
Set-TextValue 'D5' '516.11'
Set-TextValue 'E5' '  +0.17%  '
Set-TextValue 'D6' '145.52'
Set-TextValue 'E6' '  -1.31%  '
Set-TextValue 'D7' '0.999'
Set-TextValue 'E7' '  +0.37%  '
Set-TextValue 'D8' '0.570'
Set-TextValue 'E8' '  +0.04%  '
Set-TextValue 'D9' '2.654.50'
Set-TextValue 'E9' '  -0.45%  '
Set-TextValue 'D10' '6.33'
Set-TextValue 'E10' '  -3.22%  '
Set-TextValue 'E11' '  -0.76%  '
Set-TextValue 'D12' '0.336'
Set-TextValue 'E12' '  -0.83%  '
Set-TextValue 'D13' '0.127'
Set-TextValue 'E13' '  +0.52%  '
Set-TextValue 'D14' '3.108.80'
Set-TextValue 'E14' '  +0.13%  '
Set-TextValue 'D15' '59.747.92'
Set-TextValue 'E15' '  +0.76%  '
Set-TextValue 'D16' '20.87'
Set-TextValue 'E16' '  -2.18%  '
Set-TextValue 'D17' '0.0000137'
Set-TextValue 'E17' '  -1.09%  '
Set-TextValue 'D18' '2.633.16'
Set-TextValue 'E18' '  -1.15%  '
Set-TextValue 'D19' '350.16'
Set-TextValue 'E19' '  +1.31%  '
Set-TextValue 'D20' '4.49'
Set-TextValue 'E20' '  -2.37%  '
Set-TextValue 'D21' '10.31'
Set-TextValue 'E21' '  -1.83%  '
Set-TextValue 'D22' '6.19'
Set-TextValue 'E22' '  +0.09%  '
Set-TextValue 'D23' '0.999'
Set-TextValue 'E23' '  -0.06%  '
Set-TextValue 'D24' '61.98'
Set-TextValue 'E24' '  +1.40%  '
Set-TextValue 'D25' '0.415'
Set-TextValue 'E25' '  -2.53%  '
Set-TextValue 'D26' '0.164'
Set-TextValue 'E26' '  +2.25%  '
Set-TextValue 'D27' '0.997'
Set-TextValue 'E27' '  -0.39%  '
Set-TextValue 'D28' '0.0₃0804'
Set-TextValue 'E28' '  -2.54%  '
Set-TextValue 'D29' '7.09'
Set-TextValue 'E29' '  -0.52%  '
Set-TextValue 'D30' '0.998'
Set-TextValue 'E30' '  +0.10%  '
Set-TextValue 'E31' '  -3.44%  '
Set-TextValue 'B32' 'EthereumClassic'
Set-TextValue 'C32' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D32' '18.90'
Set-TextValue 'E32' '  -0.32%  '
Set-TextValue 'B33' 'PancakeSwap'
Set-TextValue 'C33' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D33' '1.57'
Set-TextValue 'E33' '  -0.20%  '
Set-TextValue 'E34' '  +0.11%  '
Set-TextValue 'D35' '0.948'
Set-TextValue 'E35' '  -10.09%  '
Set-TextValue 'D36' '4.04'
Set-TextValue 'E36' '  +0.28%  '
Set-TextValue 'D37' '1.18'
Set-TextValue 'E37' '  +2.40%  '
Set-TextValue 'D38' '0.860'
Set-TextValue 'E38' '  -1.15%  '
Set-TextValue 'D39' '36.57'
Set-TextValue 'E39' '  +0.41%  '
Set-TextValue 'E40' '  +1.09%  '
Set-TextValue 'D41' '3.65'
Set-TextValue 'E41' '  -1.21%  '
Set-TextValue 'B42' 'FirstDigitalUSD'
Set-TextValue 'C42' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D42' '1.00'
Set-TextValue 'E42' '  +0.98%  '
Set-TextValue 'D43' '0.0989'
Set-TextValue 'E43' '  -0.47%  '
Set-TextValue 'B44' 'Bittensor'
Set-TextValue 'C44' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D44' '277.02'
Set-TextValue 'E44' '  -3.23%  '
Set-TextValue 'D45' '19.69'
Set-TextValue 'E45' '  +0.75%  '
Set-TextValue 'E46' '  -3.42%  '
Set-TextValue 'D47' '2.099.97'
Set-TextValue 'E47' '  +5.72%  '
Set-TextValue 'D48' '0.0527'
Set-TextValue 'E48' '  -3.18%  '
Set-TextValue 'B49' 'WhiteBITCoin'
Set-TextValue 'C49' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D49' '10.31'
Set-TextValue 'E49' '  +0.54%  '
Set-TextValue 'D50' '0.0230'
Set-TextValue 'E50' '  -0.88%  '
Set-TextValue 'B51' 'RenderToken'
Set-TextValue 'C51' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D51' '4.70'
Set-TextValue 'E51' '  -1.36%  '
